# Commit: "Set Percentage field to 'Extended'"
#
# For every data row in the "Framework Data Model" sheet whose Component
# column (F) is "Percentage", set the Document-Support column (I) to
# "Extended". Also (re-)apply an AutoFilter over the used range, matching
# the filter state captured in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Framework Data Model")

$lastRow = 66

for ($r = 2; $r -le $lastRow; $r++) {
    $component = $ws.Cells.Item($r, 6).Value2
    if ($component -eq "Percentage") {
        $ws.Cells.Item($r, 9).Value = "Extended"
    }
}

# Re-enable the AutoFilter across the full table (A1:L66), as reflected in
# the saved workbook (adds <autoFilter> + the hidden _FilterDatabase name).
$ws.Range("A1:L66").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Framework Data Model'!`$A`$1:`$L`$66")
$filterName.Visible = $false

# Restore the selection left on the sheet after the edit.
$ws.Range("I2").Select() | Out-Null
